$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")

# Activate the "model" sheet - this becomes the new selected/active tab
# (was "choices" before the edit).
$ws.Activate()

# --- Add the new rows (74-80) describing the new "FU logic" fields ---
# Row 74: HHOID
$ws.Cells.Item(74, 1).Value = "HHOID"
$ws.Cells.Item(74, 2).Value = "integer"
$ws.Cells.Item(74, 3).Value = $false

# Row 75: BAIRRO
$ws.Cells.Item(75, 1).Value = "BAIRRO"
$ws.Cells.Item(75, 2).Value = "integer"
$ws.Cells.Item(75, 3).Value = $false
$ws.Cells.Item(75, 4).Value = "for logic"

# Row 76: HOUSEGRP
$ws.Cells.Item(76, 1).Value = "HOUSEGRP"
$ws.Cells.Item(76, 2).Value = "text"
$ws.Cells.Item(76, 3).Value = $false
$ws.Cells.Item(76, 4).Value = "for logic"

# Row 77: FAM
$ws.Cells.Item(77, 1).Value = "FAM"
$ws.Cells.Item(77, 2).Value = "integer"
$ws.Cells.Item(77, 3).Value = $false
$ws.Cells.Item(77, 4).Value = "for logic"

# Row 78: FNO
$ws.Cells.Item(78, 1).Value = "FNO"
$ws.Cells.Item(78, 2).Value = "integer"
$ws.Cells.Item(78, 3).Value = $false
$ws.Cells.Item(78, 4).Value = "fno id"

# Row 79: POID
$ws.Cells.Item(79, 1).Value = "POID"
$ws.Cells.Item(79, 2).Value = "integer"
$ws.Cells.Item(79, 3).Value = $false
$ws.Cells.Item(79, 4).Value = "x"
$ws.Cells.Item(79, 4).Value = ""

# Row 80: IDOID
$ws.Cells.Item(80, 1).Value = "IDOID"
$ws.Cells.Item(80, 2).Value = "integer"
$ws.Cells.Item(80, 3).Value = $false
$ws.Cells.Item(80, 4).Value = "x"
$ws.Cells.Item(80, 4).Value = ""

# The new cells come in unstyled (no explicit cell style index) in the
# source workbook, so strip the inherited column style back to Normal.
$ws.Range("A74:C74").Style = "Normal"
$ws.Range("A75:D78").Style = "Normal"
$ws.Range("A79:D80").Style = "Normal"

# --- Stray, empty, unstyled cells that show up in a few existing rows ---
$ws.Range("J2").Value = "x"
$ws.Range("J2").Value = ""

$ws.Range("G3").Value = "x"
$ws.Range("G3").Value = ""
$ws.Range("H3").Value = "x"
$ws.Range("H3").Value = ""
$ws.Range("I3").Value = "x"
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = "x"
$ws.Range("J3").Value = ""

$ws.Range("G4").Value = "x"
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = "x"
$ws.Range("H4").Value = ""
$ws.Range("I4").Value = "x"
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = "x"
$ws.Range("J4").Value = ""

$ws.Range("J13").Value = "x"
$ws.Range("J13").Value = ""

$ws.Range("J2").Style = "Normal"
$ws.Range("G3").Style = "Normal"
$ws.Range("H3").Style = "Normal"
$ws.Range("I3").Style = "Normal"
$ws.Range("J3").Style = "Normal"
$ws.Range("G4").Style = "Normal"
$ws.Range("H4").Style = "Normal"
$ws.Range("I4").Style = "Normal"
$ws.Range("J4").Style = "Normal"
$ws.Range("J13").Style = "Normal"

# Final selection/active cell left on the newly entered data.
$ws.Range("D77").Select()
